$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# "rewire to circuit board" -> "rewire the circuit board"
$find.Execute("rewire to circuit board", $false, $false, $false, $false, $false, `
              $true, 1, $false, "rewire the circuit board", 2)

# "of which allowed for a much cleaner wire layout" -> "which allowed for a more logical wire layout"
$find.Execute("of which allowed for a much cleaner wire layout", $false, $false, $false, $false, $false, `
              $true, 1, $false, "which allowed for a more logical wire layout", 2)
